$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2898
$ws.Range("I98").Value = 3122.5
$ws.Range("K98").Value = 3122.5
$ws.Range("M98").Value = -1624.5
$ws.Range("H122").Value = 2898
$ws.Range("I122").Value = 3122.5
$ws.Range("K122").Value = 9367.5
$ws.Range("M122").Value = -6917.5
$ws.Range("H132").Value = 1233.3334
$ws.Range("I132").Value = 1233.3334
$ws.Range("K132").Value = 3700.0002
$ws.Range("M132").Value = -1170.0002
$ws.Range("H135").Value = 770.2857
$ws.Range("I135").Value = 684.9091
$ws.Range("J135").Value = 1083.3334
$ws.Range("K135").Value = 6164.1819
$ws.Range("L135").Value = 9750.000599999999
$ws.Range("M135").Value = -3629.1819
$ws.Range("N135").Value = -14820.0006
$ws.Range("H137").Value = 2156.25
$ws.Range("I137").Value = 2010
$ws.Range("J137").Value = 2400
$ws.Range("K137").Value = 6030
$ws.Range("L137").Value = 7200
$ws.Range("M137").Value = -3480
$ws.Range("N137").Value = -12300
$ws.Range("H141").Value = 3539.8235
$ws.Range("I141").Value = 2597.4546
$ws.Range("J141").Value = 5267.5
$ws.Range("K141").Value = 7792.3638
$ws.Range("L141").Value = 15802.5
$ws.Range("M141").Value = -2612.3638
$ws.Range("N141").Value = -26162.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4715.8203
$ws.Range("I32").Value = 3307.5356
$ws.Range("K32").Value = 3307.5356
$ws.Range("M32").Value = -3020.5356
$ws.Range("H45").Value = 4738256
$ws.Range("I45").Value = 6924437.5
$ws.Range("J45").Value = 1530.1666
$ws.Range("K45").Value = 6924437.5
$ws.Range("L45").Value = 1530.1666
$ws.Range("M45").Value = -6924060.5
$ws.Range("N45").Value = -2284.1666
$ws.Range("H61").Value = 3080.1365
$ws.Range("I61").Value = 2192.5789
$ws.Range("K61").Value = 2192.5789
$ws.Range("M61").Value = -1980.5789
$ws.Range("H74").Value = 2099
$ws.Range("H77").Value = 2099
$ws.Range("H132").Value = 2942.4443
$ws.Range("I132").Value = 1998
$ws.Range("J132").Value = 3414.6667
$ws.Range("K132").Value = 5994
$ws.Range("L132").Value = 10244.0001
$ws.Range("M132").Value = -3464
$ws.Range("N132").Value = -15304.0001
$ws.Range("H136").Value = 3080.1365
$ws.Range("I136").Value = 2192.5789
$ws.Range("K136").Value = 6577.736699999999
$ws.Range("M136").Value = -4027.736699999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1724.7222
$ws.Range("I107").Value = 1406.5
$ws.Range("K107").Value = 1406.5
$ws.Range("M107").Value = 513.5
$ws.Range("H134").Value = 6817.875
$ws.Range("I134").Value = 7621.48
$ws.Range("J134").Value = 3947.8572
$ws.Range("K134").Value = 22864.44
$ws.Range("L134").Value = 11843.5716
$ws.Range("M134").Value = -20329.44
$ws.Range("N134").Value = -16913.5716
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3932.4443
$ws.Range("I31").Value = 1387
$ws.Range("J31").Value = 5968.8
$ws.Range("K31").Value = 1387
$ws.Range("L31").Value = 5968.8
$ws.Range("M31").Value = -1092
$ws.Range("N31").Value = -6558.8
$ws.Range("H34").Value = 3932.4443
$ws.Range("I34").Value = 1387
$ws.Range("J34").Value = 5968.8
$ws.Range("K34").Value = 1387
$ws.Range("L34").Value = 5968.8
$ws.Range("M34").Value = -1185
$ws.Range("N34").Value = -6372.8
$ws.Range("H59").Value = 39300
$ws.Range("J59").Value = 39300
$ws.Range("L59").Value = 39300
$ws.Range("N59").Value = -41590
$ws.Range("H99").Value = 1713.5
$ws.Range("I99").Value = 1601.1428
$ws.Range("K99").Value = 1601.1428
$ws.Range("M99").Value = -103.1428000000001
$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524
$ws.Range("H126").Value = 1713.5
$ws.Range("I126").Value = 1601.1428
$ws.Range("K126").Value = 4803.428400000001
$ws.Range("M126").Value = -2333.428400000001
$ws.Range("H134").Value = 2817.6667
$ws.Range("I134").Value = 2164.818
$ws.Range("K134").Value = 6494.454000000001
$ws.Range("M134").Value = -3959.454000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1012.3333
$ws.Range("I5").Value = 643.5
$ws.Range("K5").Value = 1930.5
$ws.Range("M5").Value = -1818.5
$ws.Range("H46").Value = 533.6667
$ws.Range("I46").Value = 750.5
$ws.Range("J46").Value = 100
$ws.Range("K46").Value = 2251.5
$ws.Range("L46").Value = 300
$ws.Range("M46").Value = -2160.5
$ws.Range("N46").Value = -482
$ws.Range("H122").Value = 1296.909
$ws.Range("I122").Value = 1021
$ws.Range("J122").Value = 1454.5714
$ws.Range("K122").Value = 9189
$ws.Range("L122").Value = 13091.1426
$ws.Range("M122").Value = -6739
$ws.Range("N122").Value = -17991.1426
$ws.Range("H131").Value = 14306769
$ws.Range("I131").Value = 71429120
$ws.Range("J131").Value = 26181.285
$ws.Range("K131").Value = 214287360
$ws.Range("L131").Value = 78543.855
$ws.Range("M131").Value = -214282320
$ws.Range("N131").Value = -88623.855
$ws.Range("H135").Value = 1012.3333
$ws.Range("I135").Value = 643.5
$ws.Range("K135").Value = 5791.5
$ws.Range("M135").Value = -3256.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 193.71428
$ws.Range("I2").Value = 212.5
$ws.Range("K2").Value = 212.5
$ws.Range("M2").Value = -99.5
$ws.Range("H64").Value = 38000
$ws.Range("J64").Value = 38000
$ws.Range("L64").Value = 38000
$ws.Range("N64").Value = -38496
$ws.Range("H67").Value = 38000
$ws.Range("J67").Value = 38000
$ws.Range("L67").Value = 38000
$ws.Range("N67").Value = -39716
$ws.Range("H97").Value = 1268.25
$ws.Range("I97").Value = 429.2
$ws.Range("J97").Value = 2666.6667
$ws.Range("K97").Value = 429.2
$ws.Range("L97").Value = 2666.6667
$ws.Range("M97").Value = 66.80000000000001
$ws.Range("N97").Value = -3658.6667
$ws.Range("H102").Value = 1893.2609
$ws.Range("J102").Value = 1682.6
$ws.Range("L102").Value = 1682.6
$ws.Range("N102").Value = -4926.6
$ws.Range("H122").Value = 2452.5
$ws.Range("I122").Value = 2254.1667
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 6762.500100000001
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -4312.500100000001
$ws.Range("N122").Value = -13150
$ws.Range("H140").Value = 10780
$ws.Range("J140").Value = 10780
$ws.Range("L140").Value = 10780
$ws.Range("N140").Value = -21140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8114.75
$ws.Range("J16").Value = 500
$ws.Range("L16").Value = 500
$ws.Range("N16").Value = -840
$ws.Range("H93").Value = 1076.4166
$ws.Range("I93").Value = 893.7
$ws.Range("J93").Value = 1990
$ws.Range("K93").Value = 893.7
$ws.Range("L93").Value = 1990
$ws.Range("M93").Value = 354.3
$ws.Range("N93").Value = -4486
$ws.Range("H132").Value = 2129.1875
$ws.Range("I132").Value = 2749
$ws.Range("J132").Value = 2040.6428
$ws.Range("K132").Value = 8247
$ws.Range("L132").Value = 6121.928400000001
$ws.Range("M132").Value = -5717
$ws.Range("N132").Value = -11181.9284
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45000
$ws.Range("J16").Value = 45000
$ws.Range("L16").Value = 45000
$ws.Range("N16").Value = -45584
$ws.Range("H100").Value = 376.77777
$ws.Range("J100").Value = 649.5
$ws.Range("L100").Value = 1299
$ws.Range("N100").Value = -2381
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
